# "modify UserMaster n Forms ctrl"
# Update the crew list: change the first crew member's details and add a
# second crew member's details on row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - existing crew member: update name / rank / email
$ws.Range("B3").Value = "amitbh"
$ws.Range("C3").Value = "Scond Officer"
$ws.Range("D3").Value = "abc@gmai.com"

# Row 4 - new crew member
$ws.Range("B4").Value = "bingshu"
$ws.Range("C4").Value = "Third Officer"
$ws.Range("D4").Value = "abc@gmai.com"

# Move the active selection, matching the saved view state
$ws.Range("F9").Select()
